$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric ("39.1" etc.) must stay stored as
# plain text (shared string), exactly like the original workbook does.
# Temporarily force a Text number format so Excel doesn't reinterpret the
# typed value as a Double, then restore the default "Normal" style so no
# extra styling is left behind on the cells.
$textCells = @("A2", "B2", "H2", "A3", "B3", "D3", "E3", "H3")

foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 updates
$ws.Range("A2").Value = "27.1"
$ws.Range("B2").Value = "186.0"
$ws.Range("H2").Value = "2025-05-04 21:39:45"

# Row 3 updates
$ws.Range("A3").Value = "40.5"
$ws.Range("B3").Value = "172.0"
$ws.Range("C3").Value = "Adelie"
$ws.Range("D3").Value = "0.91"
$ws.Range("E3").Value = "0.09"
$ws.Range("H3").Value = "2025-05-04 21:39:45"

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
